$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 12 (old rows 12+ shift down to 14+)
$ws.Range("A12:A13").EntireRow.Insert()

# Copy formatting (font, fill, number format, borders, row height) from row 11
# (an existing fully-styled data row) onto the two freshly inserted blank rows.
$ws.Range("A11:Q11").Copy()
$ws.Range("A12:Q13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 12: new item #6 - PRONTOGEST ----
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "PRONTOGEST 100MG/2ML 10 IM AMPOULE"
$ws.Range("H12").Value = "0:9"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "240.00"
$ws.Range("P12").Value = "24.0000"
$ws.Range("Q12").Value = "0:1"

# ---- Row 13: new item #7 - SPASMOFEN ----
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "SPASMOFEN 3 AMP. FOR I.M. INJ."
$ws.Range("H13").Value = "2:1"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "39.00"
$ws.Range("P13").Value = "12.8700"
$ws.Range("Q13").Value = "0:1"

# ---- Row 14 (was row 12 pre-insert, TIRATAM): renumber item to #8 ----
$ws.Range("A14").Value = 8

# ---- Row 15 (was row 13 pre-insert, سرنجات): renumber item to #9, update values ----
$ws.Range("A15").Value = 9
$ws.Range("P15").Value = "16.0000"
$ws.Range("Q15").Value = "8:0"

# ---- Row 16: totals row, update sum ----
$ws.Range("P16").Value = 535.87

# ---- Row 17: footer row, update generated timestamp ----
$ws.Range("A17").Value = "Wednesday, 27 August, 2025 10:28 AM"

# ---- Fix up merged cells for the two newly-inserted rows (Insert() shifted the
#      pre-existing merges below but did not create merges for the blank rows) ----
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()
